# The post at row 485 ("「大器晩成」") was removed from the sheet.
# Deleting its entire row shifts every subsequent row up by one
# (486->485, 487->486, ... 606->605) and shrinks the used range
# from A1:C606 down to A1:C605, matching the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(485).Delete()
